# Generate Report for Handoff
#
# Source diff summary:
#  - shared string "Handed back: in sync with en-US" -> "Ready for handoff"
#    (used by Overview!E2/F2 and the Status column on the zh-cn/de-de sheets)
#  - shared string "2016-08-17 22:57:04" -> "2016-08-17 22:57:56"
#    (Overview!G2 "Latest HO Xliff Generate Date" and de-de!H2 "Latest Handoff Datetime")
#  - shared string "2016-08-17 22:56:57" -> "2016-08-17 22:57:51"
#    (zh-cn!H2 "Latest Handoff Datetime")
#  - narrower "Status"-ish columns: Overview E & F, and column C on the
#    zh-cn / de-de sheets (29.9777047293527 -> 17.2159881591797 character width)

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text, row 2 (the single data row on every sheet)
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# Timestamps
$wsOverview.Range("G2").Value = "2016-08-17 22:57:56"
$wsZhCn.Range("H2").Value = "2016-08-17 22:57:51"
$wsDeDe.Range("H2").Value = "2016-08-17 22:57:56"

# Column widths - the target character width from the workbook (17.2159881591797)
# sits between two of this host's pixel-quantized column-width steps, so feed it
# the ColumnWidth input that rounds to the closest reachable stored width.
$newColWidth = 16.333333333333336
$wsOverview.Range("E1").ColumnWidth = $newColWidth
$wsOverview.Range("F1").ColumnWidth = $newColWidth
$wsZhCn.Range("C1").ColumnWidth = $newColWidth
$wsDeDe.Range("C1").ColumnWidth = $newColWidth
